$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 47
$ws.Range("I2").Value = 47
$ws.Range("K2").Value = 47
$ws.Range("M2").Value = 66
$ws.Range("H4").Value = 225
$ws.Range("I4").Value = 225
$ws.Range("K4").Value = 225
$ws.Range("M4").Value = -111
$ws.Range("H5").Value = 1440.125
$ws.Range("J5").Value = 1996.7273
$ws.Range("L5").Value = 1996.7273
$ws.Range("N5").Value = -2226.7273
$ws.Range("H6").Value = 119.7
$ws.Range("I6").Value = 141.14285
$ws.Range("J6").Value = 69.666664
$ws.Range("K6").Value = 423.42855
$ws.Range("L6").Value = 208.999992
$ws.Range("M6").Value = -311.42855
$ws.Range("N6").Value = -432.999992
$ws.Range("H8").Value = 10000798
$ws.Range("J8").Value = 2498.3333
$ws.Range("L8").Value = 7494.999899999999
$ws.Range("N8").Value = -7772.999899999999
$ws.Range("H9").Value = 249.5
$ws.Range("I9").Value = 249.5
$ws.Range("K9").Value = 249.5
$ws.Range("M9").Value = -80.5
$ws.Range("H12").Value = 1639.5
$ws.Range("I12").Value = 1470.8572
$ws.Range("J12").Value = 2033
$ws.Range("K12").Value = 1470.8572
$ws.Range("L12").Value = 2033
$ws.Range("M12").Value = -1300.8572
$ws.Range("N12").Value = -2373
$ws.Range("H13").Value = 5200
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5200
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5200
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -5538
$ws.Range("H16").Value = 10000
$ws.Range("J16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("N16").Value = -10460
$ws.Range("H18").Value = 27782760
$ws.Range("I18").Value = 35718324
$ws.Range("K18").Value = 35718324
$ws.Range("M18").Value = -35718040
$ws.Range("H19").Value = 6185.6113
$ws.Range("I19").Value = 1457.85
$ws.Range("J19").Value = 12095.3125
$ws.Range("K19").Value = 1457.85
$ws.Range("L19").Value = 12095.3125
$ws.Range("M19").Value = -1282.85
$ws.Range("N19").Value = -12445.3125
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H32").Value = 33335340
$ws.Range("J32").Value = 16668854
$ws.Range("L32").Value = 16668854
$ws.Range("N32").Value = -16669506
$ws.Range("H33").Value = 653.1177
$ws.Range("I33").Value = 673.8929000000001
$ws.Range("J33").Value = 556.1667
$ws.Range("K33").Value = 673.8929000000001
$ws.Range("L33").Value = 556.1667
$ws.Range("M33").Value = -444.8929000000001
$ws.Range("N33").Value = -1014.1667
$ws.Range("H38").Value = 3204.7778
$ws.Range("I38").Value = 137.6
$ws.Range("K38").Value = 412.8
$ws.Range("M38").Value = -40.79999999999995
$ws.Range("H39").Value = 430.9
$ws.Range("I39").Value = 412.22223
$ws.Range("J39").Value = 599
$ws.Range("K39").Value = 1236.66669
$ws.Range("L39").Value = 1797
$ws.Range("M39").Value = -940.66669
$ws.Range("N39").Value = -2389
$ws.Range("H40").Value = 5526.5293
$ws.Range("I40").Value = 5581.9165
$ws.Range("J40").Value = 5393.6
$ws.Range("K40").Value = 5581.9165
$ws.Range("L40").Value = 5393.6
$ws.Range("M40").Value = -5406.9165
$ws.Range("N40").Value = -5743.6
$ws.Range("H41").Value = 562.63635
$ws.Range("J41").Value = 652.8333
$ws.Range("L41").Value = 652.8333
$ws.Range("N41").Value = -1532.8333
$ws.Range("H42").Value = 1264.6666
$ws.Range("I42").Value = 1544.75
$ws.Range("J42").Value = 704.5
$ws.Range("K42").Value = 4634.25
$ws.Range("L42").Value = 2113.5
$ws.Range("M42").Value = -4404.25
$ws.Range("N42").Value = -2573.5
$ws.Range("H43").Value = 1812.5
$ws.Range("I43").Value = 1542.1111
$ws.Range("K43").Value = 1542.1111
$ws.Range("M43").Value = -1473.1111
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H51").Value = 3767.653
$ws.Range("I51").Value = 3570.6072
$ws.Range("J51").Value = 4030.3809
$ws.Range("K51").Value = 3570.6072
$ws.Range("L51").Value = 4030.3809
$ws.Range("M51").Value = -3086.6072
$ws.Range("N51").Value = -4998.3809
$ws.Range("H52").Value = 300
$ws.Range("J52").Value = 300
$ws.Range("L52").Value = 900
$ws.Range("N52").Value = -1220
$ws.Range("H129").Value = 1237.4166
$ws.Range("J129").Value = 2674.25
$ws.Range("L129").Value = 8022.75
$ws.Range("N129").Value = -18022.75
$ws.Range("H139").Value = 68200
$ws.Range("J139").Value = 68200
$ws.Range("L139").Value = 68200
$ws.Range("N139").Value = -78480
$ws.Range("H140").Value = 69696.336
$ws.Range("J140").Value = 69696.336
$ws.Range("L140").Value = 69696.336
$ws.Range("N140").Value = -80056.336

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12585727
$ws.Range("I2").Value = 14683201
$ws.Range("K2").Value = 14683201
$ws.Range("M2").Value = -14683088
$ws.Range("H116").Value = 12585727
$ws.Range("I116").Value = 14683201
$ws.Range("K116").Value = 14683201
$ws.Range("M116").Value = -14680907

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12585727
$ws.Range("I3").Value = 14683201
$ws.Range("K3").Value = 14683201
$ws.Range("M3").Value = -14683087
$ws.Range("H132").Value = 151392.4
$ws.Range("J132").Value = 151392.4
$ws.Range("L132").Value = 151392.4
$ws.Range("N132").Value = -161512.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1191.2142
$ws.Range("I16").Value = 334.625
$ws.Range("J16").Value = 2333.3333
$ws.Range("K16").Value = 334.625
$ws.Range("L16").Value = 2333.3333
$ws.Range("M16").Value = -47.625
$ws.Range("N16").Value = -2907.3333
$ws.Range("H22").Value = 13382.111
$ws.Range("I22").Value = 14992.5
$ws.Range("J22").Value = 499
$ws.Range("K22").Value = 14992.5
$ws.Range("L22").Value = 499
$ws.Range("M22").Value = -14642.5
$ws.Range("N22").Value = -1199
$ws.Range("H62").Value = 66673930
$ws.Range("I62").Value = 100007190
$ws.Range("J62").Value = 7394.8
$ws.Range("K62").Value = 100007190
$ws.Range("L62").Value = 7394.8
$ws.Range("M62").Value = -100006566
$ws.Range("N62").Value = -8642.799999999999
$ws.Range("H65").Value = 66673930
$ws.Range("I65").Value = 100007190
$ws.Range("J65").Value = 7394.8
$ws.Range("K65").Value = 500035950
$ws.Range("L65").Value = 36974
$ws.Range("M65").Value = -500032830
$ws.Range("N65").Value = -43214
$ws.Range("H113").Value = 1191.2142
$ws.Range("I113").Value = 334.625
$ws.Range("J113").Value = 2333.3333
$ws.Range("K113").Value = 334.625
$ws.Range("L113").Value = 2333.3333
$ws.Range("M113").Value = 1835.375
$ws.Range("N113").Value = -6673.3333
$ws.Range("H141").Value = 302058.12
$ws.Range("J141").Value = 323780.72
$ws.Range("L141").Value = 323780.72
$ws.Range("N141").Value = -334140.72

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 125001680
$ws.Range("I4").Value = 142857780
$ws.Range("K4").Value = 428573340
$ws.Range("M4").Value = -428573228
$ws.Range("H25").Value = 890.5
$ws.Range("I25").Value = 862.625
$ws.Range("J25").Value = 1002
$ws.Range("K25").Value = 2587.875
$ws.Range("L25").Value = 3006
$ws.Range("M25").Value = -2418.875
$ws.Range("N25").Value = -3344
$ws.Range("H26").Value = 282.18182
$ws.Range("I26").Value = 29
$ws.Range("K26").Value = 87
$ws.Range("M26").Value = 201
$ws.Range("H30").Value = 890.5
$ws.Range("I30").Value = 862.625
$ws.Range("J30").Value = 1002
$ws.Range("K30").Value = 2587.875
$ws.Range("L30").Value = 3006
$ws.Range("M30").Value = -2485.875
$ws.Range("N30").Value = -3210
$ws.Range("H55").Value = 2493.9285
$ws.Range("I55").Value = 2076.6667
$ws.Range("J55").Value = 4997.5
$ws.Range("K55").Value = 6230.000100000001
$ws.Range("L55").Value = 14992.5
$ws.Range("M55").Value = -6053.000100000001
$ws.Range("N55").Value = -15346.5
$ws.Range("H115").Value = 2296.125
$ws.Range("I115").Value = 1892.25
$ws.Range("K115").Value = 5676.75
$ws.Range("M115").Value = -4501.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 28684.459
$ws.Range("J136").Value = 28684.459
$ws.Range("L136").Value = 86053.37699999999
$ws.Range("N136").Value = -91153.37699999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 917.1667
$ws.Range("I22").Value = 917.25
$ws.Range("J22").Value = 917
$ws.Range("K22").Value = 917.25
$ws.Range("L22").Value = 917
$ws.Range("M22").Value = -622.25
$ws.Range("N22").Value = -1507
$ws.Range("H27").Value = 917.1667
$ws.Range("I27").Value = 917.25
$ws.Range("J27").Value = 917
$ws.Range("K27").Value = 917.25
$ws.Range("L27").Value = 917
$ws.Range("M27").Value = -810.25
$ws.Range("N27").Value = -1131
$ws.Range("H40").Value = 1962.619
$ws.Range("I40").Value = 1765.3529
$ws.Range("J40").Value = 2801
$ws.Range("K40").Value = 1765.3529
$ws.Range("L40").Value = 2801
$ws.Range("M40").Value = -1629.3529
$ws.Range("N40").Value = -3073

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4443.5713
$ws.Range("I132").Value = 5044.5713
$ws.Range("K132").Value = 15133.7139
$ws.Range("M132").Value = -12603.7139
